$wb = $excel.ActiveWorkbook
try { Write-Output $wb.GetType().FullName } catch { Write-Output "e1: $_" }
try { Write-Output ($wb | Get-Member -MemberType Method | Out-String) } catch { Write-Output "e2: $_" }
